# Updates cryptos list cell values (Price/Volume columns) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.885.74'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.138.46'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.03'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.47'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.136.14'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.90'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D16").Value = '3.658.50'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '63.692.31'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '3.137.74'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.07'
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.95'
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.28'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.95'
$ws.Range("E24").Value = '  -3.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.21'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.29'
$ws.Range("E28").Value = '  +7.03%  '
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '0.0₃0851'
$ws.Range("E35").Value = '  -3.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.31'
$ws.Range("E37").Value = '  -5.00%  '
$ws.Range("E38").Value = '  -2.38%  '
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '440.37'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.84'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '2.926.12'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.281'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.85'
$ws.Range("E47").Value = '  +4.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.08'
$ws.Range("E48").Value = '  +3.22%  '
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.16'
$ws.Range("E51").Value = '  -3.19%  '
